$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update Beitragsbemessungsgrenze GKV value
$ws.Range("B5").Value = 62100

# Update Jahresarbeitsentgeltgrenze GKV value
$ws.Range("B6").Value = 69300

# Update Eintragungsdatum value (kept as text)
$ws.Range("B7").Value = "01.01.2024"

# Update selection to B4
$ws.Range("B4").Select()
